# Update "想去人数" (interested-count) figures in the 北京-漫展信息 workbook.
# This mirrors a scheduled data refresh (gh-pages regeneration) that bumped
# column F on the 展览 (Exhibitions), 本地生活 (Local Life) and 全部类型
# (All Types) sheets. 演出 (Performances) is untouched.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 465
$ws.Range("F4").Value  = 500
$ws.Range("F5").Value  = 2316
$ws.Range("F6").Value  = 7
$ws.Range("F7").Value  = 63
$ws.Range("F9").Value  = 1678
$ws.Range("F10").Value = 1678
$ws.Range("F11").Value = 1387
$ws.Range("F16").Value = 749
$ws.Range("F17").Value = 182
$ws.Range("F18").Value = 127
$ws.Range("F19").Value = 7425
$ws.Range("F20").Value = 8308
$ws.Range("F23").Value = 291
$ws.Range("F25").Value = 493
$ws.Range("F26").Value = 95
$ws.Range("F28").Value = 270
$ws.Range("F29").Value = 256
$ws.Range("F33").Value = 356
$ws.Range("F34").Value = 1482
$ws.Range("F36").Value = 237
$ws.Range("F38").Value = 299
$ws.Range("F39").Value = 28
$ws.Range("F40").Value = 768
$ws.Range("F42").Value = 1372
$ws.Range("F43").Value = 363
$ws.Range("F44").Value = 264
$ws.Range("F46").Value = 93
$ws.Range("F47").Value = 202
$ws.Range("F48").Value = 183

# --- 本地生活 (sheet3) --------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2640
$ws.Range("F4").Value = 295
$ws.Range("F5").Value = 149

# --- 全部类型 (sheet4) --------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 465
$ws.Range("F5").Value  = 149
$ws.Range("F7").Value  = 500
$ws.Range("F8").Value  = 2316
$ws.Range("F9").Value  = 7
$ws.Range("F10").Value = 63
$ws.Range("F12").Value = 1678
$ws.Range("F13").Value = 1678
$ws.Range("F19").Value = 749
$ws.Range("F21").Value = 182
$ws.Range("F22").Value = 127
$ws.Range("F23").Value = 7425
$ws.Range("F24").Value = 7425
$ws.Range("F25").Value = 8308
$ws.Range("F27").Value = 293
$ws.Range("F28").Value = 95
$ws.Range("F29").Value = 270
$ws.Range("F34").Value = 237
$ws.Range("F37").Value = 299
$ws.Range("F38").Value = 28
$ws.Range("F41").Value = 768
$ws.Range("F44").Value = 1372
$ws.Range("F45").Value = 363
$ws.Range("F46").Value = 264
$ws.Range("F48").Value = 202
$ws.Range("F49").Value = 183
